$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "原文均" + "在" -> single run "原文均在", dropping the _GoBack
#    bookmark that used to sit between them (Word merges the runs and
#    removes the bookmark it used to straddle when the whole span is
#    replaced in one Find/Replace pass).
# ------------------------------------------------------------------
$find1 = $d.Content.Find
$find1.ClearFormatting()
$find1.Execute("原文均在", $true, $false, $false, $false, $false, $true, 1, $false, "原文均在", 2) | Out-Null

# ------------------------------------------------------------------
# 2) The blank paragraph right after "...zotero中。" loses the stray
#    <w:rFonts w:hint="eastAsia"/> on its paragraph-mark run
#    properties, leaving only the red color. Rewriting the (empty)
#    paragraph via InsertXML on its own whole-paragraph range
#    normalizes this cleanly (no leftover empty run).
# ------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -eq "`r") {
        $prev = $d.Paragraphs($i - 1)
        if ($prev.Range.Text -like "*zotero*") {
            $para.Range.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:rPr><w:color w:val='FF0000'/></w:rPr></w:pPr></w:p>") | Out-Null
            break
        }
    }
}

# ------------------------------------------------------------------
# 3) Split the "Cultural tree preference and its influence..." run in
#    two ("...and its in" | "fluence..."), with the _GoBack bookmark
#    now sitting between the two halves (it moved here from its old
#    spot inside "原文均在").
# ------------------------------------------------------------------
$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Execute("Cultural tree preference and its influence", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPos = $find2.Parent.Start + 35
$d.Bookmarks.Add("_GoBack", $d.Range($splitPos, $splitPos)) | Out-Null

$find3 = $d.Content.Find
$find3.ClearFormatting()
$find3.Execute("Cultural tree preference and its in", $true, $false, $false, $false, $false, $true, 1, $false, "Cultural tree preference and its in", 2) | Out-Null
